$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.732.26"
$ws.Range("E2").Value = "  +6.18%  "
$ws.Range("D3").Value = "2.395.91"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("E4").Value = "  +0.39%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "113.59"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +7.14%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "318.59"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.46%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.634"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +3.09%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.88"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.41%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0928"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.72"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.94%  "
$ws.Range("E13").Value = "  +2.41%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.81"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "2.762.22"
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").Value = "2.396.30"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "45.680.47"
$ws.Range("E18").Value = "  +6.97%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.47"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.99%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0000108"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.42"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "74.69"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.45%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.91%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "264.12"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.47%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.35"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("E26").Value = "  -0.48%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.19%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "11.31"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  +4.62%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "38.95"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.91%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "22.74"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.90%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0969"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +11.71%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "172.72"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.57%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.98"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.86%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.95%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.91"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("E37").Value = "  +4.72%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +12.88%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.08%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0362"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +12.26%  "
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "100.61"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -6.84%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.241"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.89%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.49"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +9.02%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "71.89"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "87.38"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +13.54%  "
$ws.Range("E47").Value = "  +0.14%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "115.29"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.39%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.66"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +9.25%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "9.46"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +6.20%  "
$ws.Range("D51").Value = "1.665.59"
$ws.Range("E51").Value = "  -3.22%  "
